# Fruta / hortaliza, semanal
# Re-order the weekly price records (rows 2-14) onto a new date sequence.
# Columns A,B,C,E,F,G,H,I,J,T are identical for every row, so only the
# per-record columns D,K,L,M,N,O,P,Q,R,S actually need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","K","L","M","N","O","P","Q","R","S")

# Snapshot current ("before") values for every data row so the permutation
# below can be applied safely without clobbering source data.
$snapshot = @{}
for ($r = 2; $r -le 14; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2()
    }
    $snapshot[$r] = $rowVals
}

# Mapping: new row -> old row whose values it should receive.
$mapping = @{
    2  = 8
    3  = 7
    4  = 3
    5  = 11
    6  = 12
    7  = 5
    8  = 2
    9  = 9
    10 = 13
    11 = 14
    12 = 6
    13 = 4
    14 = 10
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
